$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.531.29"
$ws.Range("E2").Value = "  +3.12%  "
$ws.Range("D3").Value = "2.436.12"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'576.10"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").Value = "'144.83"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.536"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "2.435.48"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").Value = "'0.160"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "'5.20"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "'0.351"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").Value = "'28.36"
$ws.Range("E14").Value = "  +8.66%  "
$ws.Range("E15").Value = "  +4.75%  "
$ws.Range("D16").Value = "2.878.29"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").Value = "62.474.57"
$ws.Range("E17").Value = "  +3.30%  "
$ws.Range("B18").Value = "BabyDogeCoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D18").Value = "0.0₅0102"
$ws.Range("E18").Value = "  +252.56%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "2.435.12"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'7.78"
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'10.85"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").Value = "'325.44"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  +9.46%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "'65.27"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").Value = "'625.28"
$ws.Range("E27").Value = "  +9.63%  "
$ws.Range("E28").Value = "  +13.31%  "
$ws.Range("D29").Value = "'8.44"
$ws.Range("E29").Value = "  +4.65%  "
$ws.Range("D30").Value = "0.0₃0970"
$ws.Range("E30").Value = "  +3.52%  "
$ws.Range("D31").Value = "2.558.23"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").Value = "'8.17"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").Value = "'1.40"
$ws.Range("E33").Value = "  +4.98%  "
$ws.Range("E34").Value = "  +5.98%  "
$ws.Range("D35").Value = "'1.86"
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "'4.72"
$ws.Range("E38").Value = "  +2.87%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'152.91"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").Value = "'0.371"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'5.42"
$ws.Range("E41").Value = "  +5.23%  "
$ws.Range("D42").Value = "'18.52"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "'2.71"
$ws.Range("E43").Value = "  +7.32%  "
$ws.Range("D44").Value = "'1.74"
$ws.Range("E44").Value = "  +4.03%  "
$ws.Range("D45").Value = "'42.35"
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'15.01"
$ws.Range("E47").Value = "  +28.00%  "
$ws.Range("D48").Value = "'143.46"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("D49").Value = "'3.58"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").Value = "'20.36"
$ws.Range("E50").Value = "  +5.21%  "
$ws.Range("D51").Value = "'0.600"
$ws.Range("E51").Value = "  +1.67%  "
